$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 53
$ws.Range("H53").Value = 295.60715
$ws.Range("I53").Value = 154.8125
$ws.Range("J53").Value = 483.33334
$ws.Range("K53").Value = 154.8125
$ws.Range("L53").Value = 483.33334
$ws.Range("M53").Value = 482.1875
$ws.Range("N53").Value = -1757.33334

# ALC row 80
$ws.Range("H80").Value = 1238.4117
$ws.Range("I80").Value = 615.1818
$ws.Range("J80").Value = 2381
$ws.Range("K80").Value = 1845.5454
$ws.Range("L80").Value = 7143
$ws.Range("M80").Value = -847.5454
$ws.Range("N80").Value = -9139

# ALC row 83
$ws.Range("H83").Value = 1238.4117
$ws.Range("I83").Value = 615.1818
$ws.Range("J83").Value = 2381
$ws.Range("K83").Value = 5536.6362
$ws.Range("L83").Value = 21429
$ws.Range("M83").Value = -544.6361999999999
$ws.Range("N83").Value = -31413

# ALC row 137
$ws.Range("H137").Value = 1687.3513
$ws.Range("J137").Value = 2102.682
$ws.Range("L137").Value = 6308.045999999999
$ws.Range("N137").Value = -11408.046

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 17812.266
$ws.Range("I32").Value = 17819.438
$ws.Range("K32").Value = 17819.438
$ws.Range("M32").Value = -17532.438

# ARM row 102
$ws.Range("H102").Value = 1736.6666
$ws.Range("I102").Value = 1210
$ws.Range("K102").Value = 1210
$ws.Range("M102").Value = 412

# ARM row 122
$ws.Range("H122").Value = 1295.1333
$ws.Range("I122").Value = 1296.4546
$ws.Range("J122").Value = 1291.5
$ws.Range("K122").Value = 3889.3638
$ws.Range("L122").Value = 3874.5
$ws.Range("M122").Value = -1439.3638
$ws.Range("N122").Value = -8774.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 64
$ws.Range("H64").Value = 416.7857
$ws.Range("I64").Value = 414.4
$ws.Range("J64").Value = 418.1111
$ws.Range("K64").Value = 414.4
$ws.Range("L64").Value = 418.1111
$ws.Range("M64").Value = -189.4
$ws.Range("N64").Value = -868.1111000000001

# BSM row 67
$ws.Range("H67").Value = 416.7857
$ws.Range("I67").Value = 414.4
$ws.Range("J67").Value = 418.1111
$ws.Range("K67").Value = 414.4
$ws.Range("L67").Value = 418.1111
$ws.Range("M67").Value = 365.6
$ws.Range("N67").Value = -1978.1111

# BSM row 105
$ws.Range("H105").Value = 2407.7144
$ws.Range("I105").Value = 1972.5714
$ws.Range("J105").Value = 2842.8572
$ws.Range("K105").Value = 1972.5714
$ws.Range("L105").Value = 2842.8572
$ws.Range("M105").Value = -225.5714
$ws.Range("N105").Value = -6336.8572

$ws = $wb.Worksheets.Item("CRP")
# CRP row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# CRP row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# CRP row 99
$ws.Range("H99").Value = 2391.85
$ws.Range("J99").Value = 2564.4
$ws.Range("L99").Value = 2564.4
$ws.Range("N99").Value = -5560.4

# CRP row 126
$ws.Range("H126").Value = 2391.85
$ws.Range("J126").Value = 2564.4
$ws.Range("L126").Value = 7693.200000000001
$ws.Range("N126").Value = -12633.2

$ws = $wb.Worksheets.Item("CUL")
# CUL row 16
$ws.Range("H16").Value = 199
$ws.Range("I16").Value = 199
$ws.Range("K16").Value = 597
$ws.Range("M16").Value = -424

# CUL row 39
$ws.Range("H39").Value = 4000
$ws.Range("J39").Value = 4000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588

# CUL row 62
$ws.Range("H62").Value = 2646.1538
$ws.Range("I62").Value = 1800
$ws.Range("K62").Value = 5400
$ws.Range("M62").Value = -4714

# CUL row 65
$ws.Range("H65").Value = 2646.1538
$ws.Range("I65").Value = 1800
$ws.Range("K65").Value = 16200
$ws.Range("M65").Value = -12768

# CUL row 131
$ws.Range("H131").Value = 807.61
$ws.Range("I131").Value = 400
$ws.Range("J131").Value = 811.7273
$ws.Range("K131").Value = 1200
$ws.Range("L131").Value = 2435.1819
$ws.Range("M131").Value = 3840
$ws.Range("N131").Value = -12515.1819

$ws = $wb.Worksheets.Item("GSM")
# GSM row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# GSM row 113
$ws.Range("H113").Value = 15625986
$ws.Range("I113").Value = 41667292
$ws.Range("J113").Value = 1202.6
$ws.Range("K113").Value = 41667292
$ws.Range("L113").Value = 1202.6
$ws.Range("M113").Value = -41665122
$ws.Range("N113").Value = -5542.6

# GSM row 122
$ws.Range("H122").Value = 32260550
$ws.Range("I122").Value = 76926190
$ws.Range("J122").Value = 2029.7222
$ws.Range("K122").Value = 230778570
$ws.Range("L122").Value = 6089.1666
$ws.Range("M122").Value = -230776120
$ws.Range("N122").Value = -10989.1666

# GSM row 126
$ws.Range("H126").Value = 1011663.7
$ws.Range("I126").Value = 1897.4117
$ws.Range("J126").Value = 2084540.4
$ws.Range("K126").Value = 5692.2351
$ws.Range("L126").Value = 6253621.199999999
$ws.Range("M126").Value = -3222.2351
$ws.Range("N126").Value = -6258561.199999999

# GSM row 132
$ws.Range("H132").Value = 26622.512
$ws.Range("I132").Value = 39067.445
$ws.Range("J132").Value = 2621.5715
$ws.Range("K132").Value = 117202.335
$ws.Range("L132").Value = 7864.7145
$ws.Range("M132").Value = -114672.335
$ws.Range("N132").Value = -12924.7145

$ws = $wb.Worksheets.Item("LTW")
# LTW row 76
$ws.Range("H76").Value = 9999
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 9999
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 9999
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -10675

# LTW row 79
$ws.Range("H79").Value = 9999
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 9999
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 9999
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -12339

# LTW row 100
$ws.Range("I100").Value = 66667964
$ws.Range("K100").Value = 66667964
$ws.Range("M100").Value = -66667423

# LTW row 132
$ws.Range("H132").Value = 6266.575
$ws.Range("I132").Value = 12759.765
$ws.Range("J132").Value = 1467.2609
$ws.Range("K132").Value = 38279.295
$ws.Range("L132").Value = 4401.7827
$ws.Range("M132").Value = -35749.295
$ws.Range("N132").Value = -9461.7827

# LTW row 136
$ws.Range("H136").Value = 7486.2383
$ws.Range("I136").Value = 10767.583
$ws.Range("J136").Value = 3111.111
$ws.Range("K136").Value = 32302.749
$ws.Range("L136").Value = 9333.332999999999
$ws.Range("M136").Value = -29752.749
$ws.Range("N136").Value = -14433.333

$ws = $wb.Worksheets.Item("WVR")
# WVR row 96
$ws.Range("H96").Value = 41667044
$ws.Range("I96").Value = 41667044
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 41667044
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -41665671
$ws.Range("N96").ClearContents()

# WVR row 136
$ws.Range("H136").Value = 1214.9464
$ws.Range("I136").Value = 1495.1154
$ws.Range("J136").Value = 972.13336
$ws.Range("K136").Value = 4485.3462
$ws.Range("L136").Value = 2916.40008
$ws.Range("M136").Value = -1935.3462
$ws.Range("N136").Value = -8016.40008
